$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Bet*" global constants to the new "Gacha*" constants
$ws.Range("A8").Value = "GachaEnergy"
$ws.Range("A9").Value = "Gacha1Event"
$ws.Range("A10").Value = "Gacha2Events"
$ws.Range("A11").Value = "Gacha3Events"
$ws.Range("A12").Value = "Gacha1BrokenEnergy"
$ws.Range("A13").Value = "Gacha2BrokenEnergys"

# Update the corresponding values for the renamed constants
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("D11").Value = 10
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 2

# Update the selected cell to A7
$null = $ws.Range("A7").Select()
